# Update column C ("Förändrad") date values from 2024-05-11 (serial 45423)
# to 2024-05-12 (serial 45424) for rows 2 through 28, matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2()
    if ($current -eq 45423) {
        $cell.Value = 45424
    }
}
